# "Calendrier correspondants et geko bikes"
# Insert a new calendar entry ("Dim 6 Octobre" / La Geko Bikes / UC Lutterbach VTT / VTT / gekobikes)
# between the existing "Sam 5 Octobre" (row 54) and "Sam 12 Octobre" (old row 55) entries,
# keeping the list sorted chronologically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calendrier")

# Insert a new row at position 55; this shifts the former rows 55-65 down to 56-66
# and grows the sheet's used range from F65 to F66.
$ws.Rows.Item(55).Insert()

# Populate the new row with the new event's data.
$ws.Range("A55").Value = "Dim 6 Octobre"
$ws.Range("B55").Value = "La Geko Bikes"
$ws.Range("C55").Value = "UC Lutterbach VTT"
$ws.Range("D55").Value = "VTT"
$ws.Range("E55").Value = "gekobikes"

# Match the saved view state: scrolled down so row 36 is at the top, with the
# new row's club cell (B55) selected/active.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B55").Select()
